{"js": "// Replace the arithmetic-problem text in every cell of the first (and only)\n// table in the document body with the new set of expressions, preserving\n// cell/paragraph/run formatting (Word's Table.values setter only swaps the\n// text of the existing run(s), it does not touch run/paragraph properties).\nconst newValues = [\n  [\"61+37=\", \"49-44=\", \"57-20=\", \"21-4=\", \"99-35=\"],\n  [\"37+2=\", \"73+26=\", \"0+36=\", \"79+7=\", \"71-56=\"],\n  [\"53-26=\", \"95+2=\", \"79+20=\", \"55-4=\", \"12+58=\"],\n  [\"88-54=\", \"60+19=\", \"12+20=\", \"4+89=\", \"38-26=\"],\n  [\"58+0=\", \"82-28=\", \"30-25=\", \"34+32=\", \"65-18=\"],\n  [\"71-16=\", \"47+23=\", \"7+22=\", \"86-66=\", \"91+7=\"],\n  [\"90-69=\", \"79-70=\", \"54+25=\", \"90-81=\", \"0+56=\"],\n  [\"99+0=\", \"76-52=\", \"80-56=\", \"58+29=\", \"69-30=\"],\n  [\"74-55=\", \"94-80=\", \"73+13=\", \"88-47=\", \"52-51=\"],\n  [\"28+0=\", \"92-75=\", \"5+6=\", \"39-18=\", \"99-36=\"],\n  [\"28+20=\", \"35+18=\", \"4+11=\", \"44+2=\", \"13+6=\"],\n  [\"70-55=\", \"93-7=\", \"39+29=\", \"99-21=\", \"8-0=\"],\n  [\"17+76=\", \"34+22=\", \"91-6=\", \"17+64=\", \"2+58=\"],\n  [\"88-5=\", \"31+6=\", \"69-58=\", \"27+54=\", \"96-72=\"],\n  [\"27+1=\", \"90-88=\", \"65+18=\", \"0+3=\", \"28+16=\"],\n  [\"49+29=\", \"93-27=\", \"86-33=\", \"75-55=\", \"11+6=\"],\n  [\"61-37=\", \"58+40=\", \"87-26=\", \"61-40=\", \"33-16=\"],\n  [\"2-1=\", \"34+63=\", \"5+73=\", \"10-10=\", \"42+40=\"],\n  [\"89-35=\", \"98-3=\", \"95-66=\", \"97-45=\", \"32+54=\"],\n  [\"66-21=\", \"10+64=\", \"18+35=\", \"79-79=\", \"12+4=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$newValues = @(\n  @(\"61+37=\", \"49-44=\", \"57-20=\", \"21-4=\", \"99-35=\"),\n  @(\"37+2=\", \"73+26=\", \"0+36=\", \"79+7=\", \"71-56=\"),\n  @(\"53-26=\", \"95+2=\", \"79+20=\", \"55-4=\", \"12+58=\"),\n  @(\"88-54=\", \"60+19=\", \"12+20=\", \"4+89=\", \"38-26=\"),\n  @(\"58+0=\", \"82-28=\", \"30-25=\", \"34+32=\", \"65-18=\"),\n  @(\"71-16=\", \"47+23=\", \"7+22=\", \"86-66=\", \"91+7=\"),\n  @(\"90-69=\", \"79-70=\", \"54+25=\", \"90-81=\", \"0+56=\"),\n  @(\"99+0=\", \"76-52=\", \"80-56=\", \"58+29=\", \"69-30=\"),\n  @(\"74-55=\", \"94-80=\", \"73+13=\", \"88-47=\", \"52-51=\"),\n  @(\"28+0=\", \"92-75=\", \"5+6=\", \"39-18=\", \"99-36=\"),\n  @(\"28+20=\", \"35+18=\", \"4+11=\", \"44+2=\", \"13+6=\"),\n  @(\"70-55=\", \"93-7=\", \"39+29=\", \"99-21=\", \"8-0=\"),\n  @(\"17+76=\", \"34+22=\", \"91-6=\", \"17+64=\", \"2+58=\"),\n  @(\"88-5=\", \"31+6=\", \"69-58=\", \"27+54=\", \"96-72=\"),\n  @(\"27+1=\", \"90-88=\", \"65+18=\", \"0+3=\", \"28+16=\"),\n  @(\"49+29=\", \"93-27=\", \"86-33=\", \"75-55=\", \"11+6=\"),\n  @(\"61-37=\", \"58+40=\", \"87-26=\", \"61-40=\", \"33-16=\"),\n  @(\"2-1=\", \"34+63=\", \"5+73=\", \"10-10=\", \"42+40=\"),\n  @(\"89-35=\", \"98-3=\", \"95-66=\", \"97-45=\", \"32+54=\"),\n  @(\"66-21=\", \"10+64=\", \"18+35=\", \"79-79=\", \"12+4=\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Cell($r, $c).Range.Text = $newValues[$r-1][$c-1]\n  }\n}\n"}
